# adde ReadingExcel using method and try catch
#
# Fill in the newly-read "BlockCode/Status" value (card verification flag)
# for the second data row and move the selection onto it, the way the
# reading routine leaves the cursor after it finishes populating the row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

try {
    $ws.Range("D2").Value = "V"
    $ws.Range("D3").Select()
}
catch {
    Write-Host "Failed to update CardDetails sheet: $_"
}
